$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-12 18:16:25"
$wsZh.Range("E3").Value = "2016-03-12 18:16:25"
$wsZh.Range("H2").Value = "2016-03-12 18:16:46"
$wsZh.Range("H3").Value = "2016-03-12 18:16:46"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-12 18:16:28"
$wsDe.Range("E3").Value = "2016-03-12 18:16:28"
$wsDe.Range("H2").Value = "2016-03-12 18:16:52"
$wsDe.Range("H3").Value = "2016-03-12 18:16:52"
